$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (Summary totals) and D (Total work time) for rows 2-11
$values = @(
    @{Row=2;  C=169702.5; D=449.5}
    @{Row=3;  C=89715;    D=245}
    @{Row=4;  C=57915;    D=208}
    @{Row=5;  C=69637.5;  D=210.75}
    @{Row=6;  C=91552.5;  D=210.75}
    @{Row=7;  C=115605;   D=210.25}
    @{Row=8;  C=140535;   D=210}
    @{Row=9;  C=165600;   D=210}
    @{Row=10; C=191317.5; D=210.25}
    @{Row=11; C=216337.5; D=210}
)

foreach ($item in $values) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

$wb.Save()
